# FeatherFriend / BirdDB.xlsx — add new bird records + fix SubSpec for row 2
# (commit: "unit test for register and add icons")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 2: SubSpec "12A" -> "50A"
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "50A"

# ---------------------------------------------------------------------
# 2. Row 12 (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = 54332
$ws.Cells.Item(12, 2).Value = "American Gouldian"
$ws.Cells.Item(12, 3).Value = "North America"
$ws.Cells.Item(12, 4).Value = "22A"
$ws.Cells.Item(12, 5).Value = "Male"
$ws.Cells.Item(12, 6).Value = 223
$ws.Cells.Item(12, 7).Value = 111
$ws.Cells.Item(12, 8).Value = "15/05/2023"
$ws.Cells.Item(12, 9).Value = "Red"
$ws.Cells.Item(12, 10).Value = "Purple"
$ws.Cells.Item(12, 11).Value = "Pastel"

# ---------------------------------------------------------------------
# 3. Row 13 (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = 543322
$ws.Cells.Item(13, 2).Value = "American Gouldian"
$ws.Cells.Item(13, 3).Value = "North America"
$ws.Cells.Item(13, 4).Value = "22A"
$ws.Cells.Item(13, 5).Value = "Female"
$ws.Cells.Item(13, 6).Value = 223
$ws.Cells.Item(13, 7).Value = 111
$ws.Cells.Item(13, 8).Value = "15/05/2023"
$ws.Cells.Item(13, 9).Value = "Black"
$ws.Cells.Item(13, 10).Value = "Purple"
$ws.Cells.Item(13, 11).Value = "Green"

# ---------------------------------------------------------------------
# 4. Row 14 (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = 3425
$ws.Cells.Item(14, 2).Value = "American Gouldian"
$ws.Cells.Item(14, 3).Value = "Central America"
$ws.Cells.Item(14, 4).Value = "50A"
$ws.Cells.Item(14, 5).Value = "Male"
$ws.Cells.Item(14, 6).Value = 101
$ws.Cells.Item(14, 7).Value = 1231
$ws.Cells.Item(14, 8).Value = "18/05/2023"
$ws.Cells.Item(14, 9).Value = "Red"
$ws.Cells.Item(14, 10).Value = "Purple"
$ws.Cells.Item(14, 11).Value = "Green Pastel"

# ---------------------------------------------------------------------
# 5. Row 15 (new) — H15 is a real date serial (45144) formatted as a
#    date, like H2/H3, so copy H2's format (dedupes to the same style
#    index instead of creating a new one) before writing the value.
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = 222
$ws.Cells.Item(15, 2).Value = "American Gouldian"
$ws.Cells.Item(15, 3).Value = "Central America"
$ws.Cells.Item(15, 4).Value = 444
$ws.Cells.Item(15, 5).Value = "Female"
$ws.Cells.Item(15, 6).Value = 101
$ws.Cells.Item(15, 7).Value = 1231
$ws.Range("H2").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Cells.Item(15, 8).Value = 45144
$ws.Cells.Item(15, 9).Value = "Red"
$ws.Cells.Item(15, 10).Value = "Purple"
$ws.Cells.Item(15, 11).Value = "Green Pastel"

# ---------------------------------------------------------------------
# 6. Row 16 (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = 884875
$ws.Cells.Item(16, 2).Value = "American Gouldian"
$ws.Cells.Item(16, 3).Value = "North America"
$ws.Cells.Item(16, 4).Value = "50A"
$ws.Cells.Item(16, 5).Value = "Male"
$ws.Cells.Item(16, 6).Value = 543322
$ws.Cells.Item(16, 7).Value = 111
$ws.Cells.Item(16, 8).Value = "15/05/2023"
$ws.Cells.Item(16, 9).Value = "Black"
$ws.Cells.Item(16, 10).Value = "Purple"
$ws.Cells.Item(16, 11).Value = "Green"

# ---------------------------------------------------------------------
# 7. Cosmetic: column A got very slightly narrower in the source diff
#    (19.796875 -> 19.69921875 char units). Nudge it the same direction;
#    this is the closest value the column-width model here can land on.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.0
